# "change tracing strategy and save wallet labels" - append the newly
# traced wallet hops (date + USD value at time of hop) to the bottom of
# the existing trace table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column stores its values as plain text (e.g. "2024-08-30"),
# just like every existing row. Writing a date-shaped string straight
# into .Value would get auto-recognized as a real date serial, so we
# temporarily force a text number format while assigning it, then put
# the cell's style back to the workbook default ("Normal") - leaving
# these new cells unstyled, exactly like A2:A15.
$newRows = @(
    @{ Row = 16; Date = "2024-10-05"; Usd = 0.00000091 },
    @{ Row = 17; Date = "2024-10-03"; Usd = 0.00000089 },
    @{ Row = 18; Date = "2024-01-09"; Usd = 0.00000106 }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r.Row, 2).Value = $r.Usd
}
